# Updates line-flow results (pl_mw) for the 380 kV case.
# Replaces the numeric results in columns B,C,E,F,G,H,J,L for rows 2-25
# (A2:A25 hold the time-step index 0..23) with the newly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$vals = @(1.818362043847685, 1.696117719919414, 1.621674726518108, 1.59149411452546, 1.586492061902618, 1.621267069570649, 1.776084841295074, 2.084551954139272, 2.314160155894683, 2.419265502778728, 2.459160235515697, 2.450564036930757, 2.422545793669087, 2.405395984528582, 2.307304408396703, 2.247295813553194, 2.212842183987789, 2.201187398485672, 2.253677447086886, 2.430772894938116, 2.547060423927348, 2.484945874350728, 2.250792166891131, 2.000581317326237)
$arr = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt $vals.Length; $i++) { $arr[$i,0] = $vals[$i] }
$ws.Range("B2:B25").Value = $arr

$vals = @(0.2118598411243227, 0.1955029934090078, 0.185380532939206, 0.1812356175849175, 0.1805461541182183, 0.1853247138436132, 0.2062364553657119, 0.2466169266309919, 0.2759073894301878, 0.2891519091557768, 0.2941558377081606, 0.2930786628538726, 0.2895638155023619, 0.2874093745781749, 0.2750402218524641, 0.2674316921504385, 0.2630479501447667, 0.2615623994564942, 0.2682424111100659, 0.2905965230122831, 0.3051392951345804, 0.2973836671415597, 0.2678759146526772, 0.2357595633632457)
$arr = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt $vals.Length; $i++) { $arr[$i,0] = $vals[$i] }
$ws.Range("C2:C25").Value = $arr

$vals = @(0.1743235531513694, 0.1746252053773887, 0.1748477295278779, 0.1749478198406944, 0.1749650089350858, 0.1748490412420125, 0.1744198363247929, 0.173873002580974, 0.1736494684309875, 0.1735861650505619, 0.1735676866514595, 0.1735714224087115, 0.1735845348787919, 0.1735932812491914, 0.1736543755765734, 0.1737016671074567, 0.1737324823296706, 0.1737435374133636, 0.173696258958886, 0.173580534552249, 0.1735369086294423, 0.1735572726060006, 0.173698692684674, 0.1739895289975344)
$arr = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt $vals.Length; $i++) { $arr[$i,0] = $vals[$i] }
$ws.Range("E2:E25").Value = $arr

$vals = @(2.353167595623987, 2.349591563410456, 2.34885065708221, 2.34891388011178, 2.348946412750678, 2.34885003215777, 2.351632225636038, 2.368665149630729, 2.388289208417035, 2.398772223941606, 2.402966464740416, 2.40205316006751, 2.399112781683129, 2.397340982192176, 2.38763549710356, 2.382080582305704, 2.379031945394942, 2.378024849607172, 2.38265675448821, 2.399970341387359, 2.412595001314813, 2.405736919505799, 2.382395815577141, 2.362812498845656)
$arr = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt $vals.Length; $i++) { $arr[$i,0] = $vals[$i] }
$ws.Range("F2:F25").Value = $arr

$vals = @(1.282216414801326, 1.278326469021522, 1.277021478398467, 1.276760891748253, 1.276733965674723, 1.277016867571092, 1.280649542688749, 1.296427904795053, 1.31337998580517, 1.322273969976379, 1.325813336318561, 1.325043425028525, 1.322561712592915, 1.321063955423199, 1.312822653780785, 1.308070703349216, 1.305448717285714, 1.304580016300662, 1.308565036798598, 1.323285988720585, 1.333906719724098, 1.328146302488193, 1.30834120619474, 1.291224266473506)
$arr = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt $vals.Length; $i++) { $arr[$i,0] = $vals[$i] }
$ws.Range("G2:G25").Value = $arr

$vals = @(1.197669600169419, 1.202959069354051, 1.206961449066029, 1.208781779804866, 1.209095466814091, 1.206985232676942, 1.199336527305064, 1.190345303739406, 1.187431450249392, 1.186914017706584, 1.186834791729922, 1.186846655373813, 1.186905157665478, 1.186956206910537, 1.187481570266044, 1.18801119902227, 1.188391858383028, 1.18853378593144, 1.187946946556906, 1.186884802434918, 1.186871170734435, 1.186816006905417, 1.187975757873119, 1.192131450162904)
$arr = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt $vals.Length; $i++) { $arr[$i,0] = $vals[$i] }
$ws.Range("H2:H25").Value = $arr

$vals = @(0.0661997051410177, 0.06719906058118941, 0.06784548599967266, 0.06811715843459343, 0.0681627675559362, 0.06784911647448588, 0.06653746968779473, 0.06422585402298164, 0.06268648413176159, 0.06202074684319037, 0.06177362114781459, 0.06182662268952654, 0.0620003159105611, 0.0621073561429899, 0.06273068723253061, 0.06312193095476548, 0.06335021351681025, 0.06342806390639488, 0.06307994601659761, 0.06194916291626384, 0.06123913107536971, 0.06161543170477835, 0.06309891696828185, 0.0648233030591765)
$arr = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt $vals.Length; $i++) { $arr[$i,0] = $vals[$i] }
$ws.Range("J2:J25").Value = $arr

$vals = @(0.4426690926716219, 0.4326258510333645, 0.4266081973903511, 0.424193594103059, 0.4237949299977117, 0.4265754805516337, 0.4391753788185184, 0.4650589893580559, 0.4847859918471045, 0.4939134715579172, 0.4973917563277865, 0.4966416739503785, 0.4941991941016397, 0.4927059523421065, 0.484192564575693, 0.4790090824290303, 0.4760421453952972, 0.4750400812776974, 0.4795593771890196, 0.4949160164973705, 0.5050800754466849, 0.4996437071483228, 0.4793105480624291, 0.4579316183698552)
$arr = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt $vals.Length; $i++) { $arr[$i,0] = $vals[$i] }
$ws.Range("L2:L25").Value = $arr

